# Generate Report for handoff
# Update the status of the f27a5b14-...-md file (row 3 in each sheet) from
# "Handed back: in sync with en-us" to "Ready for handoff", and refresh the
# "Latest Handoff Datetime" for the zh-cn and de-de target sheets.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# --- Overview sheet ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = $newStatus
$zhcn.Range("D3").Value = "2016-01-11 03:26:47"

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = $newStatus
$dede.Range("D3").Value = "2016-01-11 03:27:00"
